$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = "MCT-1A-Circuitos elétricos"

# Row 3
$ws.Range("B3").Value = "MCT-1A-Circuitos elétricos"
$ws.Range("C3").Value = "-"
$ws.Range("D3").Value = "MCT-3A-Máquinas Elétricas"

# Row 4
$ws.Range("B4").Value = "-"
$ws.Range("C4").Value = "-"
$ws.Range("D4").Value = "MCT-3A-Máquinas Elétricas"

# Row 6
$ws.Range("F6").Value = "-"

# Row 7
$ws.Range("B7").Value = "MCT-1A-Circuitos elétricos"
$ws.Range("E7").Value = "-"

# Row 8
$ws.Range("B8").Value = "MCT-1A-Circuitos elétricos"

# Row 18
$ws.Range("B18").Value = "[-, 'ELM-2NA-Automação Industrial', -, 'ELM-1NA-Sistemas digitais']"
$ws.Range("C18").Value = "-"

# Row 19
$ws.Range("B19").Value = "[-, 'ELM-2NA-Automação Industrial', -, 'ELM-1NA-Sistemas digitais']"
$ws.Range("C19").Value = "-"

# Row 20
$ws.Range("B20").Value = "[-, 'ELM-2NA-Automação Industrial', -, 'ELM-1NA-Sistemas digitais']"
$ws.Range("C20").Value = "-"

# Row 21
$ws.Range("B21").Value = "[-, 'ELM-2NA-Automação Industrial', -, 'ELM-1NA-Sistemas digitais']"
$ws.Range("C21").Value = "-"
